$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates scraped from the "Updated cryptos list" GitHub Actions commit.
# Column D (Price) is a text column (values like "47.102.51" / "2.40" are display
# strings, not numbers) so each D write forces text format, assigns the literal
# string, then restores the default "Normal" style so no stray number format lingers.
$updates = @(
    @{ Row=2; D="47.102.51"; E="  +5.35%  " }
    @{ Row=3; D="2.496.66"; E="  +3.06%  " }
    @{ Row=4; E="  -0.16%  " }
    @{ Row=5; D="324.44"; E="  +2.58%  " }
    @{ Row=6; D="105.79"; E="  +4.67%  " }
    @{ Row=7; D="0.524"; E="  +2.29%  " }
    @{ Row=8; D="0.999"; E="  -0.11%  " }
    @{ Row=9; E="  +2.25%  " }
    @{ Row=10; D="36.38"; E="  +3.00%  " }
    @{ Row=11; D="0.0819"; E="  +2.56%  " }
    @{ Row=12; E="  +0.83%  " }
    @{ Row=13; D="18.38"; E="  -1.62%  " }
    @{ Row=14; D="7.17"; E="  +3.89%  " }
    @{ Row=15; D="2.884.30"; E="  +3.05%  " }
    @{ Row=16; D="2.449.45"; E="  +1.88%  " }
    @{ Row=17; D="0.848"; E="  +2.07%  " }
    @{ Row=18; D="46.957.56"; E="  +5.49%  " }
    @{ Row=19; E="  +2.59%  " }
    @{ Row=20; D="6.50"; E="  +2.09%  " }
    @{ Row=21; D="0.0₃0939"; E="  +2.70%  " }
    @{ Row=22; D="70.78"; E="  +3.08%  " }
    @{ Row=23; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="251.73"; E="  +3.93%  " }
    @{ Row=24; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="2.40"; E="  +5.45%  " }
    @{ Row=25; D="2.55"; E="  +2.82%  " }
    @{ Row=26; D="26.37"; E="  +4.70%  " }
    @{ Row=28; E="  -3.49%  " }
    @{ Row=29; D="9.87"; E="  +4.13%  " }
    @{ Row=30; D="35.22"; E="  +5.29%  " }
    @{ Row=31; D="0.137"; E="  +8.47%  " }
    @{ Row=32; D="49.65"; E="  +2.51%  " }
    @{ Row=33; D="19.70"; E="  +1.35%  " }
    @{ Row=34; D="5.33"; E="  +3.38%  " }
    @{ Row=35; D="0.0772"; E="  +0.32%  " }
    @{ Row=36; E="  -0.10%  " }
    @{ Row=37; E="  +2.97%  " }
    @{ Row=38; D="1.93"; E="  +2.32%  " }
    @{ Row=39; D="2.97"; E="  +4.28%  " }
    @{ Row=40; D="123.09"; E="  -1.47%  " }
    @{ Row=41; E="  +1.95%  " }
    @{ Row=42; D="2.21"; E="  +0.73%  " }
    @{ Row=43; D="20.89"; E="  -0.92%  " }
    @{ Row=44; D="0.0295"; E="  +1.99%  " }
    @{ Row=45; D="1.982.04"; E="  +2.23%  " }
    @{ Row=46; D="2.99"; E="  +1.96%  " }
    @{ Row=47; E="  -0.70%  " }
    @{ Row=48; D="1.81"; E="  +4.55%  " }
    @{ Row=49; B="THORChain"; C="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D="5.41"; E="  +17.98%  " }
    @{ Row=50; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="9.04"; E="  -2.01%  " }
    @{ Row=51; D="80.19"; E="  +6.05%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B$($u.Row)").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$($u.Row)").Value = $u.C }
    if ($u.ContainsKey("D")) {
        $d = $ws.Range("D$($u.Row)")
        $d.NumberFormat = "@"
        $d.Value = $u.D
        $d.Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Range("E$($u.Row)").Value = $u.E }
}
